# Update "想去人数" (interested-count) figures for both the "展览" and
# "全部类型" worksheets, which carry duplicate data tables.
$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 820
    "F6"  = 12306
    "F7"  = 51
    "F10" = 437
    "F11" = 1130
    "F12" = 903
    "F13" = 13600
    "F14" = 13772
    "F19" = 1029
    "F22" = 4862
    "F23" = 212
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
